$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "第74期 第四代寵物"
$ws.Range("A24").Value = "9/5"
$ws.Range("B24").Value = "10/31"

$ws.Range("A25").Value = "9/12"
$ws.Range("B25").Value = "11/7"
$ws.Range("C25").Value = "第75期 9轉技能(新增新的通用技能池-用風暴卡抽 五層可以選40張通用池或自選九轉專屬技能 二選一)"

$ws.Range("C25").Select()
